# NatMI LR-pairs output (Amh -> Acvr1, OldD0) was regenerated against an
# updated TPM expression matrix ("update scripts wuth new tpm").
#
# Sheet layout: rows 2-10 are the 3x3 (Sending cluster x Target cluster)
# grid over clusters ECs / FAPs / MuSCs (each cluster has 3 cells total):
#   row  2: ECs   -> ECs      row  5: FAPs  -> ECs      row  8: MuSCs -> ECs
#   row  3: ECs   -> FAPs     row  6: FAPs  -> FAPs     row  9: MuSCs -> FAPs
#   row  4: ECs   -> MuSCs    row  7: FAPs  -> MuSCs    row 10: MuSCs -> MuSCs
#
# Columns E-J (ligand/sending-cluster stats) depend only on the sending
# cluster; columns K-P (receptor/target-cluster stats) depend only on the
# target cluster; columns Q-T (edge weight/specificity) multiplicatively
# combine the two sides. With the refreshed TPM matrix, only the ECs
# cluster's underlying ligand (Amh) and receptor (Acvr1) numbers moved
# (one more cell now detects Amh, and both clusters' mean expression
# shifted), which ripples into every column except A-D, K and L.
#
# Values below are the recomputed NatMI outputs for every affected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1153286666666667
$ws.Range("H2").Value = 0.345986
$ws.Range("I2").Value = 0.1212009326543938
$ws.Range("J2").Value = 0.1212009326543938
$ws.Range("M2").Value = 4.621579
$ws.Range("N2").Value = 13.864737
$ws.Range("O2").Value = 0.1778708528171788
$ws.Range("P2").Value = 0.1778708528171788
$ws.Range("Q2").Value = 0.5330005439646667
$ws.Range("R2").Value = 4.797004895682
$ws.Range("S2").Value = 0.02155811325347448
$ws.Range("T2").Value = 0.02155811325347448

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1153286666666667
$ws.Range("H3").Value = 0.345986
$ws.Range("I3").Value = 0.1212009326543938
$ws.Range("J3").Value = 0.1212009326543938
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5971062807549863
$ws.Range("P3").Value = 0.5971062807549863
$ws.Range("Q3").Value = 1.789264330869556
$ws.Range("R3").Value = 16.103378977826
$ws.Range("S3").Value = 0.07236983812130067
$ws.Range("T3").Value = 0.07236983812130066

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1153286666666667
$ws.Range("H4").Value = 0.345986
$ws.Range("I4").Value = 0.1212009326543938
$ws.Range("J4").Value = 0.1212009326543938
$ws.Range("O4").Value = 0.2250228664278349
$ws.Range("P4").Value = 0.2250228664278349
$ws.Range("Q4").Value = 0.6742943450875556
$ws.Range("R4").Value = 6.068649105787999
$ws.Range("S4").Value = 0.02727298127961868
$ws.Range("T4").Value = 0.02727298127961868

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.5471739422864045
$ws.Range("J5").Value = 0.5471739422864045
$ws.Range("M5").Value = 4.621579
$ws.Range("N5").Value = 13.864737
$ws.Range("O5").Value = 0.1778708528171788
$ws.Range("P5").Value = 0.1778708528171788
$ws.Range("Q5").Value = 2.406285186877
$ws.Range("R5").Value = 21.656566681893
$ws.Range("S5").Value = 0.09732629575382054
$ws.Range("T5").Value = 0.09732629575382055

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.5471739422864045
$ws.Range("J6").Value = 0.5471739422864045
$ws.Range("N6").Value = 46.543441
$ws.Range("O6").Value = 0.5971062807549863
$ws.Range("P6").Value = 0.5971062807549863
$ws.Range("R6").Value = 72.70034286414901
$ws.Range("S6").Value = 0.3267209976046785
$ws.Range("T6").Value = 0.3267209976046785

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.5471739422864045
$ws.Range("J7").Value = 0.5471739422864045
$ws.Range("O7").Value = 0.2250228664278349
$ws.Range("P7").Value = 0.2250228664278349
$ws.Range("S7").Value = 0.1231266489279055
$ws.Range("T7").Value = 0.1231266489279055

# Row 8 (MuSCs -> ECs)
$ws.Range("H8").Value = 0.9466730000000001
$ws.Range("I8").Value = 0.3316251250592017
$ws.Range("J8").Value = 0.3316251250592017
$ws.Range("M8").Value = 4.621579
$ws.Range("N8").Value = 13.864737
$ws.Range("O8").Value = 0.1778708528171788
$ws.Range("P8").Value = 0.1778708528171788
$ws.Range("Q8").Value = 1.458374685555667
$ws.Range("R8").Value = 13.125372170001
$ws.Range("S8").Value = 0.05898644380988378
$ws.Range("T8").Value = 0.05898644380988379

# Row 9 (MuSCs -> FAPs)
$ws.Range("H9").Value = 0.9466730000000001
$ws.Range("I9").Value = 0.3316251250592017
$ws.Range("J9").Value = 0.3316251250592017
$ws.Range("N9").Value = 46.543441
$ws.Range("O9").Value = 0.5971062807549863
$ws.Range("P9").Value = 0.5971062807549863
$ws.Range("Q9").Value = 4.895713213532556
$ws.Range("R9").Value = 44.06141892179301
$ws.Range("S9").Value = 0.1980154450290071
$ws.Range("T9").Value = 0.1980154450290071

# Row 10 (MuSCs -> MuSCs)
$ws.Range("H10").Value = 0.9466730000000001
$ws.Range("I10").Value = 0.3316251250592017
$ws.Range("J10").Value = 0.3316251250592017
$ws.Range("O10").Value = 0.2250228664278349
$ws.Range("P10").Value = 0.2250228664278349
$ws.Range("S10").Value = 0.07462323622031081
$ws.Range("T10").Value = 0.07462323622031081
